# Generate Report for handback
# This script updates the localization-status workbook to reflect that the
# "7e297ee6-...md" file (and its dependent "be8126c9-...md") have now been
# handed back / are in sync with en-US, instead of merely "Ready for handoff".
#
# It touches three worksheets:
#   - Overview: summary status cells for the zh-cn / de-de columns
#   - zh-cn:    per-language detail row, including new "Latest Target File"
#               and "Latest Handback File" links + refreshed handback datetime
#   - de-de:    same as zh-cn but for the de-de language

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: rows 3 & 4 (7e297ee6.md and be8126c9.md) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both languages.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack
$wsOverview.Range("B4").Value = $statusHandedBack
$wsOverview.Range("C4").Value = $statusHandedBack

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 3 & 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B3").Value = $statusHandedBack
$wsZh.Range("B4").Value = $statusHandedBack

# New "Latest Target File" (E) / "Latest Handback File" (F) links, mirroring
# the style already used on the existing hyperlink cells in the row (A3/C3).
$wsZh.Range("E3").Value = "7e297ee6-02d8-45aa-a999-6defc4b6135f.md"
$wsZh.Range("E3").Font.Underline = 2
$wsZh.Range("E3").Font.Color = 15631086
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/88f7d3ab2c8d71a22cd4e269454cf9d75ceb9e6c/e2e/7e297ee6-02d8-45aa-a999-6defc4b6135f.md", "", "", "7e297ee6-02d8-45aa-a999-6defc4b6135f.md")

$wsZh.Range("F3").Value = "7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.zh-cn.xlf"
$wsZh.Range("F3").Font.Underline = 2
$wsZh.Range("F3").Font.Color = 15631086
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/88f7d3ab2c8d71a22cd4e269454cf9d75ceb9e6c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.zh-cn.xlf", "", "", "7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.zh-cn.xlf")

$wsZh.Range("G3").Value = "2016-01-26 05:42:25"

$wsZh.Range("E4").Value = "7e297ee6-02d8-45aa-a999-6defc4b6135f.md"
$wsZh.Range("E4").Font.Underline = 2
$wsZh.Range("E4").Font.Color = 15631086
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/88f7d3ab2c8d71a22cd4e269454cf9d75ceb9e6c/e2e/7e297ee6-02d8-45aa-a999-6defc4b6135f.md", "", "", "7e297ee6-02d8-45aa-a999-6defc4b6135f.md")

$wsZh.Range("F4").Value = "7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.zh-cn.xlf"
$wsZh.Range("F4").Font.Underline = 2
$wsZh.Range("F4").Font.Color = 15631086
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/88f7d3ab2c8d71a22cd4e269454cf9d75ceb9e6c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.zh-cn.xlf", "", "", "7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.zh-cn.xlf")

$wsZh.Range("G4").Value = "2016-01-26 05:42:25"

# ---------------------------------------------------------------------------
# de-de sheet: rows 3 & 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B3").Value = $statusHandedBack
$wsDe.Range("B4").Value = $statusHandedBack

$wsDe.Range("E3").Value = "7e297ee6-02d8-45aa-a999-6defc4b6135f.md"
$wsDe.Range("E3").Font.Underline = 2
$wsDe.Range("E3").Font.Color = 15631086
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/be2760260187ce5b742b5c375e19363e53a974e4/e2e/7e297ee6-02d8-45aa-a999-6defc4b6135f.md", "", "", "7e297ee6-02d8-45aa-a999-6defc4b6135f.md")

$wsDe.Range("F3").Value = "7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.de-de.xlf"
$wsDe.Range("F3").Font.Underline = 2
$wsDe.Range("F3").Font.Color = 15631086
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/be2760260187ce5b742b5c375e19363e53a974e4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.de-de.xlf", "", "", "7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.de-de.xlf")

$wsDe.Range("G3").Value = "2016-01-26 05:42:49"

$wsDe.Range("E4").Value = "7e297ee6-02d8-45aa-a999-6defc4b6135f.md"
$wsDe.Range("E4").Font.Underline = 2
$wsDe.Range("E4").Font.Color = 15631086
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/be2760260187ce5b742b5c375e19363e53a974e4/e2e/7e297ee6-02d8-45aa-a999-6defc4b6135f.md", "", "", "7e297ee6-02d8-45aa-a999-6defc4b6135f.md")

$wsDe.Range("F4").Value = "7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.de-de.xlf"
$wsDe.Range("F4").Font.Underline = 2
$wsDe.Range("F4").Font.Color = 15631086
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/be2760260187ce5b742b5c375e19363e53a974e4/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.de-de.xlf", "", "", "7e297ee6-02d8-45aa-a999-6defc4b6135f.3685c224ed984c6ac98f3ec04d5a211e9fb0e530.de-de.xlf")

$wsDe.Range("G4").Value = "2016-01-26 05:42:49"
